$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header cell format (bold, bordered, centered) from H1 into I1:J1
# so the new header cells match the existing header styling.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New per-game data for columns I (I0) and J (IF), rows 2-62
$IVals = @(9,9,9,9,9,9,9,9,9,9,6,9,7,8,9,8,10,9,9,8,8,8,8,8,8,7,10,7,7,8,8,7,7,8,11,5,7,9,8,8,7,8,5,7,8,9,8,8,6,7,8,11,9,8,6,9,6,7,8,6,4)
$JVals = @(9,9,9,9,9,9,9,9,9,9,6,9,7,8,9,9,10,9,9,8,8,8,8,8,8,7,10,7,7,8,8,8,7,8,11,5,7,9,8,8,7,8,6,7,8,9,8,8,6,7,8,11,9,8,6,9,6,7,8,6,4)

for ($i = 0; $i -lt $IVals.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 9).Value = $IVals[$i]
    $ws.Cells.Item($row, 10).Value = $JVals[$i]
}

Write-Output "Added I0/IF columns for $($IVals.Length) data rows"
